$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the affected range so numeric-looking strings
# (e.g. "584.27") are stored as text, matching the inlineStr source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.523.88"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.262.63"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "584.27"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "182.05"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "0.425"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").Value = "3.841.19"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D14").Value = "28.63"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "68.552.02"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "3.285.56"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "394.91"
$ws.Range("E20").Value = "  +4.79%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +4.45%  "
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "5.69"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "22.91"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "164.59"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").Value = "0.828"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "26.26"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").Value = "41.33"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "344.31"
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("D46").Value = "2.613.30"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "0.0282"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").Value = "31.61"
$ws.Range("E51").Value = "  -0.19%  "

# Restore default (unstyled) cell style now that values are set as text.
$ws.Range("D2:E51").Style = "Normal"
